{"js": "const replacements = [\n  [\"2025-11-09 Sunday\", \"2025-11-10 Monday\"],\n  [\"112\u00d77=784\", \"957\u00d75=4785\"],\n  [\"151\u00d79=1359\", \"752\u00d78=6016\"],\n  [\"153\u00d77=1071\", \"559\u00d72=1118\"],\n  [\"474\u00d77=3318\", \"114\u00d74=456\"],\n  [\"784\u00d77=5488\", \"845\u00d76=5070\"],\n  [\"202\u00d72=404\", \"582\u00d75=2910\"],\n  [\"105\u00d76=630\", \"172\u00d72=344\"],\n  [\"525\u00d76=3150\", \"967\u00d77=6769\"],\n  [\"302\u00d75=1510\", \"195\u00d72=390\"],\n  [\"979\u00d78=7832\", \"367\u00d74=1468\"],\n  [\"562\u00d79=5058\", \"287\u00d78=2296\"],\n  [\"311\u00d76=1866\", \"194\u00d74=776\"],\n  [\"329\u00d79=2961\", \"939\u00d79=8451\"],\n  [\"720\u00d75=3600\", \"363\u00d77=2541\"],\n  [\"622\u00d79=5598\", \"324\u00d79=2916\"],\n  [\"755\u00d76=4530\", \"136\u00d72=272\"],\n  [\"613\u00d74=2452\", \"800\u00d78=6400\"],\n  [\"580\u00d78=4640\", \"944\u00d74=3776\"],\n  [\"163\u00d76=978\", \"629\u00d72=1258\"],\n  [\"166\u00d73=498\", \"199\u00d74=796\"],\n  [\"627\u00d78=5016\", \"259\u00d76=1554\"],\n  [\"314\u00d76=1884\", \"506\u00d72=1012\"],\n  [\"594\u00d73=1782\", \"335\u00d72=670\"],\n  [\"455\u00d77=3185\", \"575\u00d79=5175\"],\n  [\"528\u00d77=3696\", \"653\u00d77=4571\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-11-09 Sunday', '2025-11-10 Monday'),\n    @('112\u00d77=784', '957\u00d75=4785'),\n    @('151\u00d79=1359', '752\u00d78=6016'),\n    @('153\u00d77=1071', '559\u00d72=1118'),\n    @('474\u00d77=3318', '114\u00d74=456'),\n    @('784\u00d77=5488', '845\u00d76=5070'),\n    @('202\u00d72=404', '582\u00d75=2910'),\n    @('105\u00d76=630', '172\u00d72=344'),\n    @('525\u00d76=3150', '967\u00d77=6769'),\n    @('302\u00d75=1510', '195\u00d72=390'),\n    @('979\u00d78=7832', '367\u00d74=1468'),\n    @('562\u00d79=5058', '287\u00d78=2296'),\n    @('311\u00d76=1866', '194\u00d74=776'),\n    @('329\u00d79=2961', '939\u00d79=8451'),\n    @('720\u00d75=3600', '363\u00d77=2541'),\n    @('622\u00d79=5598', '324\u00d79=2916'),\n    @('755\u00d76=4530', '136\u00d72=272'),\n    @('613\u00d74=2452', '800\u00d78=6400'),\n    @('580\u00d78=4640', '944\u00d74=3776'),\n    @('163\u00d76=978', '629\u00d72=1258'),\n    @('166\u00d73=498', '199\u00d74=796'),\n    @('627\u00d78=5016', '259\u00d76=1554'),\n    @('314\u00d76=1884', '506\u00d72=1012'),\n    @('594\u00d73=1782', '335\u00d72=670'),\n    @('455\u00d77=3185', '575\u00d79=5175'),\n    @('528\u00d77=3696', '653\u00d77=4571'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
